$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '79.345.25'
$ws.Cells.Item(2, 5).Value = '  +3.86%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.198.97'
$ws.Cells.Item(3, 5).Value = '  +7.38%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''207.52'
$ws.Cells.Item(5, 5).Value = '  +3.50%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''631.14'
$ws.Cells.Item(6, 5).Value = '  -0.01%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +13.50%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.579'
$ws.Cells.Item(9, 5).Value = '  +5.63%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '3.195.37'
$ws.Cells.Item(10, 5).Value = '  +7.30%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.591'
$ws.Cells.Item(11, 5).Value = '  +37.45%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +1.55%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''5.43'
$ws.Cells.Item(13, 5).Value = '  +8.82%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.780.20'
$ws.Cells.Item(14, 5).Value = '  +7.31%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''0.0000229'
$ws.Cells.Item(15, 5).Value = '  +21.56%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''31.82'
$ws.Cells.Item(16, 5).Value = '  +9.84%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '79.058.38'
$ws.Cells.Item(17, 5).Value = '  +3.63%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.175.51'
$ws.Cells.Item(18, 5).Value = '  +6.84%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''14.44'
$ws.Cells.Item(19, 5).Value = '  +7.51%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''9.46'
$ws.Cells.Item(20, 5).Value = '  +6.17%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''430.50'
$ws.Cells.Item(21, 5).Value = '  +15.66%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''2.85'
$ws.Cells.Item(22, 5).Value = '  +25.90%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''4.93'
$ws.Cells.Item(23, 5).Value = '  +14.82%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''6.82'
$ws.Cells.Item(24, 5).Value = '  +6.26%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +10.07%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'WrappedeETH'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(26, 4).Value = '3.352.24'
$ws.Cells.Item(26, 5).Value = '  +6.97%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Litecoin'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(27, 4).Value = '''77.06'
$ws.Cells.Item(27, 5).Value = '  +5.97%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Aptos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(28, 4).Value = '''11.05'
$ws.Cells.Item(28, 5).Value = '  +13.35%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''1.00'
$ws.Cells.Item(29, 5).Value = '  +0.26%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''0.0000115'
$ws.Cells.Item(30, 5).Value = '  +8.67%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''1.00'
$ws.Cells.Item(31, 5).Value = '  -0.08%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''8.96'
$ws.Cells.Item(32, 5).Value = '  +8.84%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''1.48'
$ws.Cells.Item(33, 5).Value = '  +7.10%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''522.59'
$ws.Cells.Item(34, 5).Value = '  +3.19%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +1.58%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''0.137'
$ws.Cells.Item(36, 5).Value = '  +23.61%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''22.76'
$ws.Cells.Item(37, 5).Value = '  +12.31%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''0.998'
$ws.Cells.Item(38, 5).Value = '  -0.21%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''0.401'
$ws.Cells.Item(39, 5).Value = '  +5.48%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''164.26'
$ws.Cells.Item(40, 5).Value = '  +0.13%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''197.77'
$ws.Cells.Item(41, 5).Value = '  +6.71%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''0.111'
$ws.Cells.Item(42, 5).Value = '  +5.59%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +0.08%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.52%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''5.45'
$ws.Cells.Item(45, 5).Value = '  +10.18%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''0.802'
$ws.Cells.Item(46, 5).Value = '  +14.29%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +9.20%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''1.31'
$ws.Cells.Item(48, 5).Value = '  +6.25%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'OKB'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(49, 4).Value = '''42.94'
$ws.Cells.Item(49, 5).Value = '  +1.72%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'dogwifhat'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(50, 4).Value = '''2.56'
$ws.Cells.Item(50, 5).Value = '  +10.69%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''0.629'
$ws.Cells.Item(51, 5).Value = '  +6.96%  '
